# "Add practice to chapter 4 slide"
#
# Slide 30 ("Vocabulary: Learning curve") is reworked into a new
# "Practice: Code smell" slide: title + body text are replaced, the two
# wikipedia learning-curve graphs are removed, and the two remaining text
# boxes are repositioned/resized. The slide's speaker notes (which held the
# old "Use learning curve theory to..." talking points) are cleared too.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(30)

# ---------------------------------------------------------------------
# Title textbox ("object 2"): "Vocabulary: Learning curve" -> "Practice: Code smell"
# ---------------------------------------------------------------------
$title = $s.Shapes.Item(1)

$title.Left = 18.0
$title.Top = 1143476.0 / 12700.0
$title.Width = 930.0
$title.Height = 615553.0 / 12700.0

$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Practice: Code smell"

# "Practice" keeps the accent-blue color; ": Code smell" reverts to the
# default (black) text color instead of inheriting the blue from run 1.
$titleRange.Characters(1, 8).Font.Color.RGB = 15773696   # 00B0F0
$titleRange.Characters(9, 13).Font.Color.RGB = 0          # 000000

# ---------------------------------------------------------------------
# Body textbox ("object 3"): new practice instructions
# ---------------------------------------------------------------------
$body = $s.Shapes.Item(2)

$body.Left = 892986.0 / 12700.0
$body.Top = 2367695.0 / 12700.0
$body.Width = 10534219.0 / 12700.0
$body.Height = 3939540.0 / 12700.0

$bodyRange = $body.TextFrame.TextRange
$cr = [char]13
$lines = @(
    "Find a project on Github and review it's code",
    "",
    "Definition of Done",
    "List 5 naming issues",
    "List 3 function issues",
    "List 3 comment issues",
    "",
    "Create a pull request from your edit suggestion"
)
$bodyRange.Text = [string]::Join($cr, $lines)

# Highlight the "Definition of Done" paragraph in yellow like the rest of
# the deck's callouts.
$bodyRange.Paragraphs(3, 1).Font.Color.RGB = 65535   # FFFF00

# ---------------------------------------------------------------------
# Drop the two learning-curve graph pictures - no longer relevant.
# ---------------------------------------------------------------------
$s.Shapes.Item("Picture 2").Delete()
$s.Shapes.Item("Picture 4").Delete()

# ---------------------------------------------------------------------
# Clear the now-stale speaker notes for this slide.
# ---------------------------------------------------------------------
$notes = $s.NotesPage
$notesBody = $notes.Shapes.Item(2)
$notesBody.TextFrame.TextRange.Text = ""
